$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @("Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity")
for ($i = 0; $i -lt $row1.Length; $i++) { $ws.Cells.Item(1, $i+1).Value = $row1[$i] }

$row2 = @("ECs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 20.35014666666666, 61.05043999999999, 0.7972172161636617, 0.7972172161636618, 3, 1, 3.825035, 11.475105, 0.03111562857396839, 0.03111562857396839, 77.84002325513332, 700.5602092961999, 0.02480591479092156, 0.02480591479092157)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i+1).Value = $row2[$i] }

$row3 = @("ECs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 20.35014666666666, 61.05043999999999, 0.7972172161636617, 0.7972172161636618, 3, 1, 57.89762366666667, 173.692871, 0.4709815605157605, 0.4709815605157605, 1178.225133268138, 10604.02619941324, 0.3754746085387917, 0.3754746085387918)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i+1).Value = $row3[$i] }

$row4 = @("ECs", "Pdgfd", "Pdgfrb", "MuSCs", 3, 1, 20.35014666666666, 61.05043999999999, 0.7972172161636617, 0.7972172161636618, 3, 1, 61.10114166666667, 183.303425, 0.4970413158429724, 0.4970413158429724, 1243.417194417444, 11190.754749757, 0.3962498941346578, 0.3962498941346578)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i+1).Value = $row4[$i] }

$row5 = @("ECs", "Pdgfd", "Pdgfrb", "Resolving-Mac", 3, 1, 20.35014666666666, 61.05043999999999, 0.7972172161636617, 0.7972172161636618, 3, 1, 0.1059033333333333, 0.31771, 0.0008614950672987739, 0.0008614950672987739, 2.155148365822222, 19.3963352924, 0.0006867986992906549, 0.000686798699290655)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i+1).Value = $row5[$i] }

$row6 = @("FAPs", "Pdgfd", "Pdgfrb", "ECs", 3, 1, 4.304558, 12.913674, 0.1686311062905535, 0.1686311062905535, 3, 1, 3.825035, 11.475105, 0.03111562857396839, 0.03111562857396839, 16.46508500953, 148.18576508577, 0.005247062869354247, 0.005247062869354247)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i+1).Value = $row6[$i] }

$row7 = @("FAPs", "Pdgfd", "Pdgfrb", "FAPs", 3, 1, 4.304558, 12.913674, 0.1686311062905535, 0.1686311062905535, 3, 1, 57.89762366666667, 173.692871, 0.4709815605157605, 0.4709815605157605, 249.2236791353394, 2243.013112218054, 0.07942214159222395, 0.07942214159222395)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i+1).Value = $row7[$i] }

$row8 = @("FAPs", "Pdgfd", "Pdgfrb", "MuSCs", 3, 1, 4.304558, 12.913674, 0.1686311062905535, 0.1686311062905535, 3, 1, 61.10114166666667, 183.303425, 0.4970413158429724, 0.4970413158429724, 263.0134081703833, 2367.12067353345, 0.08381662696271285, 0.08381662696271285)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i+1).Value = $row8[$i] }

$row9 = @("FAPs", "Pdgfd", "Pdgfrb", "Resolving-Mac", 3, 1, 4.304558, 12.913674, 0.1686311062905535, 0.1686311062905535, 3, 1, 0.1059033333333333, 0.31771, 0.0008614950672987739, 0.0008614950672987739, 0.4558670407266667, 4.10280336654, 0.0001452748662624471, 0.0001452748662624471)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i+1).Value = $row9[$i] }

$row10 = @("MuSCs", "Pdgfd", "Pdgfrb", "ECs", 2, 0.6666666666666666, 0.8288763333333332, 2.486629, 0.03247123933933695, 0.03247123933933695, 3, 1, 3.825035, 11.475105, 0.03111562857396839, 0.03111562857396839, 3.170480985671666, 28.53432887104499, 0.001010363022619239, 0.001010363022619239)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i+1).Value = $row10[$i] }

$row11 = @("MuSCs", "Pdgfd", "Pdgfrb", "FAPs", 2, 0.6666666666666666, 0.8288763333333332, 2.486629, 0.03247123933933695, 0.03247123933933695, 3, 1, 57.89762366666667, 173.692871, 0.4709815605157605, 0.4709815605157605, 47.98997001353988, 431.909730121859, 0.01529335497592166, 0.01529335497592166)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i+1).Value = $row11[$i] }

$row12 = @("MuSCs", "Pdgfd", "Pdgfrb", "MuSCs", 2, 0.6666666666666666, 0.8288763333333332, 2.486629, 0.03247123933933695, 0.03247123933933695, 3, 1, 61.10114166666667, 183.303425, 0.4970413158429724, 0.4970413158429724, 50.64529026714722, 455.807612404325, 0.01613954752827612, 0.01613954752827612)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i+1).Value = $row12[$i] }

$row13 = @("MuSCs", "Pdgfd", "Pdgfrb", "Resolving-Mac", 2, 0.6666666666666666, 0.8288763333333332, 2.486629, 0.03247123933933695, 0.03247123933933695, 3, 1, 0.1059033333333333, 0.31771, 0.0008614950672987739, 0.0008614950672987739, 0.0877807666211111, 0.7900268995899999, 0.00002797381251991668, 0.00002797381251991668)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i+1).Value = $row13[$i] }

$row14 = @("Resolving-Mac", "Pdgfd", "Pdgfrb", "ECs", 1, 0.3333333333333333, 0.04289566666666667, 0.128687, 0.001680438206447867, 0.001680438206447868, 3, 1, 3.825035, 11.475105, 0.03111562857396839, 0.03111562857396839, 0.1640774263483333, 1.476696837135, 0.00005228789107333745, 0.00005228789107333746)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i+1).Value = $row14[$i] }

$row15 = @("Resolving-Mac", "Pdgfd", "Pdgfrb", "FAPs", 1, 0.3333333333333333, 0.04289566666666667, 0.128687, 0.001680438206447867, 0.001680438206447868, 3, 1, 57.89762366666667, 173.692871, 0.4709815605157605, 0.4709815605157605, 2.483557165597444, 22.352014490377, 0.0007914554088231222, 0.0007914554088231223)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i+1).Value = $row15[$i] }

$row16 = @("Resolving-Mac", "Pdgfd", "Pdgfrb", "MuSCs", 1, 0.3333333333333333, 0.04289566666666667, 0.128687, 0.001680438206447867, 0.001680438206447868, 3, 1, 61.10114166666667, 183.303425, 0.4970413158429724, 0.4970413158429724, 2.620974205886111, 23.588767852975, 0.0008352472173256525, 0.0008352472173256526)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i+1).Value = $row16[$i] }

$row17 = @("Resolving-Mac", "Pdgfd", "Pdgfrb", "Resolving-Mac", 1, 0.3333333333333333, 0.04289566666666667, 0.128687, 0.001680438206447867, 0.001680438206447868, 3, 1, 0.1059033333333333, 0.31771, 0.0008614950672987739, 0.0008614950672987739, 0.004542794085555556, 0.04088514677, 0.000001447689225755236, 0.000001447689225755237)
for ($i = 0; $i -lt $row17.Length; $i++) { $ws.Cells.Item(17, $i+1).Value = $row17[$i] }
